$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.441.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3684"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07500"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.976"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06758"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.345"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.441.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.392"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.628"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.063"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.745.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.017"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.760"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08332"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("E37").Value = "  +0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2286"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.339"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06494"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6202"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "

$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.773"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.049"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07329"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "

Write-Host "Cryptos list updated"